# "Generate Report for Handoff"
#
# The localization-status report moves from "In Translation" to
# "Ready for handoff": the Overview sheet's per-language status columns
# (zh-cn / de-de) and each language sheet's own Status column get the new
# label, the Latest HO Xliff Generate Date / Latest Handoff Datetime
# timestamps are refreshed, and the Status columns are widened so the
# longer "Ready for handoff" text still fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ----------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # Overview!zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # Overview!de-de status
$wsZhCn.Range("C2").Value = "Ready for handoff"       # zh-cn!Status
$wsDeDe.Range("C2").Value = "Ready for handoff"       # de-de!Status

# --- Refreshed handoff timestamps -----------------------------------------
$wsOverview.Range("G2").Value = "2016-08-22 00:50:51" # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value = "2016-08-22 00:50:51"     # de-de!Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-22 00:50:47"     # zh-cn!Latest Handoff Datetime

# --- Widen the Status columns so "Ready for handoff" fits ------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.35
$wsOverview.Columns.Item(6).ColumnWidth = 16.35
$wsZhCn.Columns.Item(3).ColumnWidth = 16.35
$wsDeDe.Columns.Item(3).ColumnWidth = 16.35
